{"js": "// 1) \"r_s : Series resistance\" table cell \u2014 collapse the split\n//    \" : S\" / \"eries\" / \" resistance\" runs (with the proofErr spell-check\n//    markers in between) into one clean run, same visible text.\nconst seriesResults = context.document.body.search(\" : Series resistance\", { matchCase: true });\nseriesResults.load(\"items\");\nawait context.sync();\n\nif (seriesResults.items.length > 0) {\n  seriesResults.items[0].insertText(\" : Series resistance\", \"Replace\");\n  await context.sync();\n}\n\n// 2) \"-FFT\" bullet line \u2014 append the additional comma separated\n//    techniques that were added after it.\nconst fftResults = context.document.body.search(\"-FFT\", { matchCase: true });\nfftResults.load(\"items\");\nawait context.sync();\n\nif (fftResults.items.length > 0) {\n  fftResults.items[0].insertText(\n    \", shifted discrete Fourier translations (SDFT), discrete wavelet transform (DWT), frequency weighting functions\",\n    \"After\"\n  );\n  await context.sync();\n}\n", "ps1": "# 1) \"r_s : Series resistance\" table cell \u2014 collapse the split\n#    \" : S\" / \"eries\" / \" resistance\" runs (with the proofErr spell-check\n#    markers in between) into one clean run, same visible text.\n$d = $word.ActiveDocument\n\n$seriesRange = $d.Content\n$found = $seriesRange.Find.Execute(\" : Series resistance\", $false, $false, $false, $false, $false, $true, 1, $false, \" : Series resistance\", 2)\nif (-not $found) {\n    Write-Output \"WARNING: series resistance text not found\"\n}\n\n# 2) \"-FFT\" bullet line \u2014 append the additional comma separated\n#    techniques that were added after it.\n$fftRange = $d.Content\n$found2 = $fftRange.Find.Execute(\"-FFT\")\nif ($found2) {\n    $fftRange.Collapse(0)\n    $fftRange.InsertAfter(\", shifted discrete Fourier translations (SDFT), discrete wavelet transform (DWT), frequency weighting functions\")\n}\n"}
